$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.586.68"
$ws.Range("E2").Value = '  -0.82%  '

$ws.Range("D3").Value = "'1.751.93"
$ws.Range("E3").Value = '  +0.08%  '

$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").Value = "'324.42"
$ws.Range("E5").Value = '  +1.03%  '

$ws.Range("D6").Value = "'1.004"
$ws.Range("E6").Value = '  +0.05%  '

$ws.Range("D7").Value = "'0.4584"
$ws.Range("E7").Value = '  +8.02%  '

$ws.Range("D8").Value = "'0.3573"
$ws.Range("E8").Value = '  -1.94%  '

$ws.Range("D9").Value = "'0.07486"
$ws.Range("E9").Value = '  +1.10%  '

$ws.Range("D10").Value = "'42.11"
$ws.Range("E10").Value = '  -2.01%  '

$ws.Range("D11").Value = "'1.093"
$ws.Range("E11").Value = '  +0.48%  '

$ws.Range("E12").Value = '  -0.19%  '

$ws.Range("E13").Value = '  +0.03%  '

$ws.Range("D14").Value = "'6.000"
$ws.Range("E14").Value = '  -0.99%  '

$ws.Range("D15").Value = "'7.077"
$ws.Range("E15").Value = '  -2.82%  '

$ws.Range("D16").Value = "'1.755.68"
$ws.Range("E16").Value = '  -1.91%  '

$ws.Range("D17").Value = "'92.23"
$ws.Range("E17").Value = '  +0.77%  '

$ws.Range("D18").Value = "'0.00001064"
$ws.Range("E18").Value = '  +0.76%  '

$ws.Range("D19").Value = "'0.06422"
$ws.Range("E19").Value = '  +0.86%  '

$ws.Range("E20").Value = '  -0.11%  '

$ws.Range("E21").Value = '  -1.42%  '

$ws.Range("D22").Value = "'5.804"
$ws.Range("E22").Value = '  -2.20%  '

$ws.Range("D23").Value = "'27.634.39"
$ws.Range("E23").Value = '  -0.81%  '

$ws.Range("D24").Value = "'11.22"
$ws.Range("E24").Value = '  -0.29%  '

$ws.Range("D25").Value = "'2.117"
$ws.Range("E25").Value = '  +1.63%  '

$ws.Range("D26").Value = "'164.51"
$ws.Range("E26").Value = '  +4.78%  '

$ws.Range("D27").Value = "'20.28"
$ws.Range("E27").Value = '  +0.88%  '

$ws.Range("D28").Value = "'1.955.08"
$ws.Range("E28").Value = '  -1.46%  '

$ws.Range("D29").Value = "'2.077"
$ws.Range("E29").Value = '  -2.85%  '

$ws.Range("D30").Value = "'126.16"
$ws.Range("E30").Value = '  +1.27%  '

$ws.Range("D31").Value = "'1.059"
$ws.Range("E31").Value = '  -6.41%  '

$ws.Range("D32").Value = "'0.09176"
$ws.Range("E32").Value = '  +3.23%  '

$ws.Range("D33").Value = "'3.668"
$ws.Range("E33").Value = '  +0.47%  '

$ws.Range("D34").Value = "'5.525"
$ws.Range("E34").Value = '  -0.71%  '

$ws.Range("D35").Value = "'11.87"
$ws.Range("E35").Value = '  -3.21%  '

$ws.Range("D36").Value = "'0.02292"
$ws.Range("E36").Value = '  +0.31%  '

$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").Value = "'0.06042"
$ws.Range("E37").Value = '  +0.85%  '

$ws.Range("B38").Value = 'Algorand'
$ws.Range("C38").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D38").Value = "'0.2095"
$ws.Range("E38").Value = '  -0.11%  '

$ws.Range("D39").Value = "'4.965"
$ws.Range("E39").Value = '  -0.09%  '

$ws.Range("D40").Value = "'0.6318"
$ws.Range("E40").Value = '  +0.07%  '

$ws.Range("D41").Value = "'1.210"
$ws.Range("E41").Value = '  +3.03%  '

$ws.Range("E42").Value = '  -1.35%  '

$ws.Range("D43").Value = "'7.771"
$ws.Range("E43").Value = '  -0.35%  '

$ws.Range("D44").Value = "'13.18"
$ws.Range("E44").Value = '  -1.45%  '

$ws.Range("D45").Value = "'0.5896"
$ws.Range("E45").Value = '  +0.37%  '

$ws.Range("D46").Value = "'3.717"
$ws.Range("E46").Value = '  +0.89%  '

$ws.Range("D47").Value = "'123.29"
$ws.Range("E47").Value = '  +0.91%  '

$ws.Range("D48").Value = "'1.938"
$ws.Range("E48").Value = '  -2.03%  '

$ws.Range("D49").Value = "'1.140"
$ws.Range("E49").Value = '  -2.87%  '

$ws.Range("D50").Value = "'0.06857"
$ws.Range("E50").Value = '  +0.53%  '

$ws.Range("D51").Value = "'71.96"
$ws.Range("E51").Value = '  -2.52%  '
